# Created experiment order generation script
# Re-names each task-order sheet (keeping its tab position / sheetId / r:id
# unchanged) and rewrites its "task_order" rows with freshly generated
# stimulus-file names, growing or shrinking the used range as needed.

$wb = $excel.ActiveWorkbook

function Set-TaskOrderSheet {
    param($Worksheet, $NewName, $Values)

    $Worksheet.Name = $NewName

    # Determine how many data rows currently exist (row 1 is the header).
    $usedRange = $Worksheet.UsedRange
    $currentLastRow = $usedRange.Rows.Count

    $newCount = $Values.Count
    $newLastRow = $newCount + 1

    # If we need more rows than currently exist, stamp the style of the
    # last existing "index" cell (column A) down over the new rows first,
    # so the new cells pick up the same bold/border/alignment formatting.
    if ($newLastRow -gt $currentLastRow) {
        $srcAddr = "A" + $currentLastRow
        $dstAddr = "A" + ($currentLastRow + 1) + ":A" + $newLastRow
        $Worksheet.Range($srcAddr).Copy($Worksheet.Range($dstAddr))
    }

    # Write the new index (0-based) and file-name values.
    for ($i = 0; $i -lt $newCount; $i++) {
        $row = $i + 2
        $Worksheet.Range("A" + $row).Value = $i
        $Worksheet.Range("B" + $row).Value = $Values[$i]
    }

    # If the sheet used to have more rows than we need now, delete the
    # surplus rows entirely so the dimension shrinks back down.
    if ($newLastRow -lt $currentLastRow) {
        $rowSpec = ($newLastRow + 1).ToString() + ":" + $currentLastRow.ToString()
        $Worksheet.Rows($rowSpec).Delete()
    }
}

# --- Sheet at tab position 1 (was GNG_TO) -> becomes NB_TO ---------------
$values1 = @(
    "OB-16515890080283446.csv",
    "OB-1651589008078596.csv",
    "ZB-match_4-1651589007814657.csv",
    "OB-16515890079890049.csv",
    "TB-16515890095806267.csv",
    "ZB-match_3-16515890078458734.csv",
    "TB-16515890089957044.csv",
    "ZB-match_0-16515890078927443.csv",
    "TB-16515890094846842.csv"
)
Set-TaskOrderSheet $wb.Worksheets.Item(1) "NB_TO-1651589009611877" $values1

# --- Sheet at tab position 2 (was NB_TO) -> becomes TOL_TO ----------------
$values2 = @(
    "MM_stims-16515890096274955.csv",
    "ZM_stims-1651589009611877.csv",
    "MM_stims-16515890096431239.csv",
    "ZM_stims-16515890096274955.csv",
    "MM_stims-1651589009658752.csv",
    "ZM_stims-16515890096431239.csv"
)
Set-TaskOrderSheet $wb.Worksheets.Item(2) "TOL_TO-1651589009658752" $values2

# --- Sheet at tab position 3 (was RS_TO) -> becomes GNG_TO ----------------
$values3 = @(
    "go_stims-1651589009658752.csv",
    "GNG_stims-16515890096743402.csv",
    "go_stims-16515890096743402.csv",
    "GNG_stims-16515890096899672.csv"
)
Set-TaskOrderSheet $wb.Worksheets.Item(3) "GNG_TO-16515890096899672" $values3

# --- Sheet at tab position 4 (was TOL_TO) -> becomes RS_TO ----------------
$values4 = @(
    "eyes closed",
    "eyes open"
)
Set-TaskOrderSheet $wb.Worksheets.Item(4) "RS_TO-16515890096899672" $values4

# --- Sheet at tab position 5 (was vSAT_TO) -> stays vSAT_TO --------------
$values5 = @(
    "vSAT_stims-16515890097212162.csv",
    "SAT_stims-16515890097055924.csv",
    "vSAT_stims-1651589009736843.csv",
    "SAT_stims-16515890096899672.csv"
)
Set-TaskOrderSheet $wb.Worksheets.Item(5) "vSAT_TO-16515890097525496" $values5
